# DADA_week7_lab2.docx edit script
$d = $word.ActiveDocument

# --- 1. Merge the "assigned 'points' to each URL..." sentence into one run,
#        append the new closing sentence about higher point value. ---
$old1 = "Our team approached the lab by first coming up with a set of rules that would indicate a possibly malicious website.  Then we assigned " + [char]0x2018 + "points" + [char]0x2019 + " to each URL based on our rules."
$new1 = "Our team approached the lab by first coming up with a set of rules that would indicate a possibly malicious website.  Then we assigned a weighted point value to each URL based on our rules.  A higher point value indicates a greater degree of potential malicious behavior."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- 2. "org." -> "org" ---
$d.Content.Find.Execute("org.", $true, $false, $false, $false, $false, $true, 1, $false, "org", 2) | Out-Null

# --- 3. ".pdf files" -> ".pdf or .exe files" ---
$d.Content.Find.Execute("file extensions (like .pdf files).", $true, $false, $false, $false, $false, $true, 1, $false, "file extensions (like .pdf or .exe files).", 2) | Out-Null

# --- 4. Point threshold sentence replaced entirely ---
$old4 = "A website reaching a certain point threshold (yet to be determined) will be considered malicious."
$new4 = "Our program then classifies any website with a score greater than four as malicious."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# --- 5. "We also considered avoiding" -> "We also attempted avoid" ---
$d.Content.Find.Execute("We also considered avoiding false positives by awarding ", $true, $false, $false, $false, $false, $true, 1, $false, "We also attempted avoid false positives by awarding ", 2) | Out-Null

# --- 6. "has characteristics" -> "had characteristics" ---
$old6 = "meaning that if a URL has characteristics indicating legitimacy (like being located in the same region it was registered) the URL could gain points back."
$new6 = "meaning that if a URL had characteristics indicating legitimacy (like being located in the same region it was registered) the URL could gain points back."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# --- 7. Append two new paragraphs at the end of the document ---

# 7a. Blank paragraph
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Last

# 7b. Results paragraph with mixed formatting runs
$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()

$resultsPara = $d.Paragraphs.Last
$resultsPara.Range.Text = "Our final run against the classify.json file found 953 legitimate links and 1,071 malicious links."

# Italicize "classify.json"
$r = $d.Content
$r.Find.Execute("classify.json", $true) | Out-Null
$r.Italic = $true
$r.Font.Name = "Baskerville"

# Bold-ish SemiBold font for the numbers
$r2 = $d.Content
$r2.Find.Execute("953", $true) | Out-Null
$r2.Font.Name = "Baskerville SemiBold"

$r3 = $d.Content
$r3.Find.Execute("1,071", $true) | Out-Null
$r3.Font.Name = "Baskerville SemiBold"

Write-Host "Edit script completed"
